# Scheduled runner: refresh market-board price snapshots (currentAveragePrice*,
# LevePrice*, LeveProfit*) per-sheet/per-row from the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 177.125
$ws.Range("I18").Value = 177.125
$ws.Range("K18").Value = 177.125
$ws.Range("M18").Value = 106.875
$ws.Range("H40").Value = 2597
$ws.Range("I40").Value = 2466.1538
$ws.Range("J40").Value = 2786
$ws.Range("K40").Value = 2466.1538
$ws.Range("L40").Value = 2786
$ws.Range("M40").Value = -2291.1538
$ws.Range("N40").Value = -3136
$ws.Range("H80").Value = 672.85187
$ws.Range("I80").Value = 225
$ws.Range("J80").Value = 1088.7142
$ws.Range("K80").Value = 675
$ws.Range("L80").Value = 3266.1426
$ws.Range("M80").Value = 323
$ws.Range("N80").Value = -5262.142599999999
$ws.Range("H83").Value = 672.85187
$ws.Range("I83").Value = 225
$ws.Range("J83").Value = 1088.7142
$ws.Range("K83").Value = 2025
$ws.Range("L83").Value = 9798.427799999999
$ws.Range("M83").Value = 2967
$ws.Range("N83").Value = -19782.4278
$ws.Range("H87").Value = 16454.303
$ws.Range("J87").Value = 16454.303
$ws.Range("L87").Value = 16454.303
$ws.Range("N87").Value = -18950.303
$ws.Range("H90").Value = 16454.303
$ws.Range("J90").Value = 16454.303
$ws.Range("L90").Value = 49362.909
$ws.Range("N90").Value = -61842.909
$ws.Range("H101").Value = 457.92856
$ws.Range("I101").Value = 453.41666
$ws.Range("K101").Value = 1360.24998
$ws.Range("M101").Value = 261.7500199999999
$ws.Range("H129").Value = 975.8421
$ws.Range("I129").Value = 465.66666
$ws.Range("J129").Value = 1071.5
$ws.Range("K129").Value = 1396.99998
$ws.Range("L129").Value = 3214.5
$ws.Range("M129").Value = 3603.00002
$ws.Range("N129").Value = -13214.5
$ws.Range("H138").Value = 1719.6
$ws.Range("I138").Value = 1378.5834
$ws.Range("J138").Value = 2034.3846
$ws.Range("K138").Value = 4135.7502
$ws.Range("L138").Value = 6103.1538
$ws.Range("M138").Value = 1004.2498
$ws.Range("N138").Value = -16383.1538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H32").Value = 4539.44
$ws.Range("I32").Value = 3876.8171
$ws.Range("J32").Value = 13342.857
$ws.Range("K32").Value = 3876.8171
$ws.Range("L32").Value = 13342.857
$ws.Range("M32").Value = -3589.8171
$ws.Range("N32").Value = -13916.857
$ws.Range("H54").Value = 11000
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 13000
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 13000
$ws.Range("M54").Value = -4231
$ws.Range("N54").Value = -14538
$ws.Range("H61").Value = 2563.0386
$ws.Range("I61").Value = 2487
$ws.Range("J61").Value = 2706.6667
$ws.Range("K61").Value = 2487
$ws.Range("L61").Value = 2706.6667
$ws.Range("M61").Value = -2275
$ws.Range("N61").Value = -3130.6667
$ws.Range("H74").Value = 1420.7906
$ws.Range("I74").Value = 985.7646999999999
$ws.Range("J74").Value = 3064.2222
$ws.Range("K74").Value = 985.7646999999999
$ws.Range("L74").Value = 3064.2222
$ws.Range("M74").Value = -111.7646999999999
$ws.Range("N74").Value = -4812.2222
$ws.Range("H77").Value = 1420.7906
$ws.Range("I77").Value = 985.7646999999999
$ws.Range("J77").Value = 3064.2222
$ws.Range("K77").Value = 4928.8235
$ws.Range("L77").Value = 15321.111
$ws.Range("M77").Value = -560.8234999999995
$ws.Range("N77").Value = -24057.111
$ws.Range("H136").Value = 2563.0386
$ws.Range("I136").Value = 2487
$ws.Range("J136").Value = 2706.6667
$ws.Range("K136").Value = 7461
$ws.Range("L136").Value = 8120.000100000001
$ws.Range("M136").Value = -4911
$ws.Range("N136").Value = -13220.0001
$ws.Range("H139").Value = 45811.07
$ws.Range("J139").Value = 45811.07
$ws.Range("L139").Value = 45811.07
$ws.Range("N139").Value = -56091.07

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 34999
$ws.Range("J2").Value = 34999
$ws.Range("L2").Value = 34999
$ws.Range("N2").Value = -35225
$ws.Range("H133").Value = 50124.125
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H134").Value = 1496
$ws.Range("I134").Value = 1226.9333
$ws.Range("J134").Value = 1944.4445
$ws.Range("K134").Value = 3680.7999
$ws.Range("L134").Value = 5833.333500000001
$ws.Range("M134").Value = -1145.7999
$ws.Range("N134").Value = -10903.3335
$ws.Range("H140").Value = 54412.855
$ws.Range("J140").Value = 54412.855
$ws.Range("L140").Value = 54412.855
$ws.Range("N140").Value = -64772.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1169.6
$ws.Range("I22").Value = 1287.1111
$ws.Range("J22").Value = 112
$ws.Range("K22").Value = 1287.1111
$ws.Range("L22").Value = 112
$ws.Range("M22").Value = -937.1111000000001
$ws.Range("N22").Value = -812
$ws.Range("H132").Value = 737386
$ws.Range("I132").Value = 1088327.2
$ws.Range("J132").Value = 3599.818
$ws.Range("K132").Value = 3264981.6
$ws.Range("L132").Value = 10799.454
$ws.Range("M132").Value = -3262451.6
$ws.Range("N132").Value = -15859.454
$ws.Range("H134").Value = 2835.9092
$ws.Range("I134").Value = 3192.25
$ws.Range("J134").Value = 1885.6666
$ws.Range("K134").Value = 9576.75
$ws.Range("L134").Value = 5656.9998
$ws.Range("M134").Value = -7041.75
$ws.Range("N134").Value = -10726.9998
$ws.Range("H138").Value = 38185.383
$ws.Range("J138").Value = 38185.383
$ws.Range("L138").Value = 38185.383
$ws.Range("N138").Value = -48465.383
$ws.Range("H140").Value = 65521.145
$ws.Range("J140").Value = 65521.145
$ws.Range("L140").Value = 65521.145
$ws.Range("N140").Value = -75881.14499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4915.5557
$ws.Range("I3").Value = 620
$ws.Range("J3").Value = 6142.857
$ws.Range("K3").Value = 1860
$ws.Range("L3").Value = 18428.571
$ws.Range("M3").Value = -1748
$ws.Range("N3").Value = -18652.571
$ws.Range("H80").Value = 3724.3333
$ws.Range("I80").Value = 3951
$ws.Range("J80").Value = 3679
$ws.Range("K80").Value = 11853
$ws.Range("L80").Value = 11037
$ws.Range("M80").Value = -10917
$ws.Range("N80").Value = -12909
$ws.Range("H83").Value = 3724.3333
$ws.Range("I83").Value = 3951
$ws.Range("J83").Value = 3679
$ws.Range("K83").Value = 35559
$ws.Range("L83").Value = 33111
$ws.Range("M83").Value = -30879
$ws.Range("N83").Value = -42471
$ws.Range("H109").Value = 2477.5
$ws.Range("I109").Value = 521
$ws.Range("J109").Value = 5999.2
$ws.Range("K109").Value = 1563
$ws.Range("L109").Value = 17997.6
$ws.Range("M109").Value = -523
$ws.Range("N109").Value = -20077.6
$ws.Range("H114").Value = 174.77777
$ws.Range("I114").Value = 152.875
$ws.Range("J114").Value = 350
$ws.Range("K114").Value = 458.625
$ws.Range("L114").Value = 1050
$ws.Range("M114").Value = 2795.375
$ws.Range("N114").Value = -7558
$ws.Range("H131").Value = 854.4299999999999
$ws.Range("J131").Value = 873.6391599999999
$ws.Range("L131").Value = 2620.91748
$ws.Range("N131").Value = -12700.91748
$ws.Range("H136").Value = 1979.5834
$ws.Range("I136").Value = 1965
$ws.Range("K136").Value = 5895
$ws.Range("M136").Value = -795
$ws.Range("H140").Value = 2009.6666
$ws.Range("I140").Value = 1549.4762
$ws.Range("J140").Value = 2815
$ws.Range("K140").Value = 4648.4286
$ws.Range("L140").Value = 8445
$ws.Range("M140").Value = 531.5713999999998
$ws.Range("N140").Value = -18805

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1043.45
$ws.Range("I2").Value = 874.25
$ws.Range("J2").Value = 1297.25
$ws.Range("K2").Value = 874.25
$ws.Range("L2").Value = 1297.25
$ws.Range("M2").Value = -761.25
$ws.Range("N2").Value = -1523.25
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H135").Value = 39878.89
$ws.Range("J135").Value = 39878.89
$ws.Range("L135").Value = 39878.89
$ws.Range("N135").Value = -50018.89

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 28459
$ws.Range("I40").Value = 38841.43
$ws.Range("J40").Value = 4233.3335
$ws.Range("K40").Value = 38841.43
$ws.Range("L40").Value = 4233.3335
$ws.Range("M40").Value = -38705.43
$ws.Range("N40").Value = -4505.3335
$ws.Range("H46").Value = 1284.8334
$ws.Range("I46").Value = 903.5
$ws.Range("J46").Value = 1589.9
$ws.Range("K46").Value = 903.5
$ws.Range("L46").Value = 1589.9
$ws.Range("M46").Value = -715.5
$ws.Range("N46").Value = -1965.9
$ws.Range("H99").Value = 31900
$ws.Range("J99").Value = 31900
$ws.Range("L99").Value = 31900
$ws.Range("N99").Value = -37890
$ws.Range("H106").Value = 15456.667
$ws.Range("J106").Value = 15456.667
$ws.Range("L106").Value = 15456.667
$ws.Range("N106").Value = -17980.667
$ws.Range("H133").Value = 44679.668
$ws.Range("J133").Value = 44679.668
$ws.Range("L133").Value = 44679.668
$ws.Range("N133").Value = -49739.668
